$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Ajout de tâches dans le planner & modification fichier excel"
$ws.Range("B3").Value = 45369
$ws.Range("C3").Value = 1

$ws.Range("A4").Value = "Ajout en base de données d'une colonne qui spécifie le rôle par défaut ou non. Modification de l'API pour empêcher la modification d'un rôle par défaut. Ajout d'un cadenas rouge bloqué et vert débloqué pour différencier les deux"
$ws.Range("B4").Value = 45369
$ws.Range("C4").Value = 1

$ws.Range("A5").Value = "Maintenance du planner"
$ws.Range("B5").Value = 45369
$ws.Range("C5").Value = 0.5

$ws.Range("B2").Copy() | Out-Null
$ws.Range("B3:B5").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

$ws.Range("C6").Select()
